# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the status that used to
# read "Handed back: in sync with en-US" becomes "Ready for handoff", and
# the associated timestamps are bumped a couple of minutes later. The
# "Status"/"zh-cn"/"de-de" columns on the per-language sheets also get
# narrower now that the status text is shorter.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$overview.Range("F2").Value = "Ready for handoff"   # de-de status
$zhcn.Range("C2").Value = "Ready for handoff"       # Status column
$dede.Range("C2").Value = "Ready for handoff"       # Status column

# --- Timestamps bumped to reflect the new handoff generation ---
$overview.Range("G2").Value = "2016-10-18 04:53:54" # Latest HO Xliff Generate Date
$dede.Range("H2").Value = "2016-10-18 04:53:54"     # Latest Handoff Datetime
$zhcn.Range("H2").Value = "2016-10-18 04:53:39"      # Latest Handoff Datetime

# --- Narrower Status/zh-cn/de-de columns to fit the shorter text ---
$overview.Columns.Item(5).ColumnWidth = 16.33  # column E (zh-cn)
$overview.Columns.Item(6).ColumnWidth = 16.33  # column F (de-de)
$zhcn.Columns.Item(3).ColumnWidth = 16.33      # column C (Status)
$dede.Columns.Item(3).ColumnWidth = 16.33      # column C (Status)
